$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Convert H:MM:SS-style times to total-minutes:SS (the "TODO in minuten umrechnen?" task)
$ws.Range("F38").Value = "65:34"
$ws.Range("E39").Value = "66:15"
$ws.Range("F39").Value = "78:19"
$ws.Range("E40").Value = "78:49"
$ws.Range("F40").Value = "90:20"
$ws.Range("E41").Value = "91:24"
$ws.Range("F41").Value = "94:35"
$ws.Range("F49").Value = "64:39"
$ws.Range("F64").Value = "60:02"
$ws.Range("E65").Value = "60:43"
$ws.Range("F65").Value = "62:16"

# 2) Remove the now-obsolete TODO comment cell
$ws.Range("H38").ClearContents()

# 3) Fill in the last participants (14, 15, 16) results
$ws.Range("C106").Value = 4
$ws.Range("D106").NumberFormat = "@"
$ws.Range("D106").Value = "0"
$ws.Range("D106").NumberFormat = "0"
$ws.Range("E106").Value = "2:48"
$ws.Range("F106").Value = "5:43"

$ws.Range("C107").Value = 1
$ws.Range("D107").Value = 1
$ws.Range("E107").Value = "6:56"
$ws.Range("F107").Value = "14:41"

$ws.Range("C108").Value = 1
$ws.Range("D108").NumberFormat = "@"
$ws.Range("D108").Value = "0"
$ws.Range("D108").NumberFormat = "0"
$ws.Range("E108").Value = "16:07"
$ws.Range("F108").Value = "17:56"

$ws.Range("C109").Value = 0
$ws.Range("D109").Value = 2
$ws.Range("E109").Value = "18:35"
$ws.Range("F109").Value = "22:23"

$ws.Range("C110").Value = 8
$ws.Range("D110").NumberFormat = "@"
$ws.Range("D110").Value = "0"
$ws.Range("D110").NumberFormat = "0"
$ws.Range("E110").Value = "25:48"
$ws.Range("F110").Value = "27:21"

$ws.Range("C111").Value = 2
$ws.Range("D111").NumberFormat = "@"
$ws.Range("D111").Value = "0"
$ws.Range("D111").NumberFormat = "0"
$ws.Range("E111").Value = "28:17"
$ws.Range("F111").Value = "29:22"

$ws.Range("C112").Value = 3
$ws.Range("D112").NumberFormat = "@"
$ws.Range("D112").Value = "0"
$ws.Range("D112").NumberFormat = "0"
$ws.Range("E112").Value = "30:14"
$ws.Range("F112").Value = "32:45"

$ws.Range("C113").Value = 0
$ws.Range("D113").Value = 1
$ws.Range("E113").Value = "33:30"
$ws.Range("F113").Value = "40:58"

$ws.Range("C114").Value = 2
$ws.Range("D114").NumberFormat = "@"
$ws.Range("D114").Value = "0"
$ws.Range("D114").NumberFormat = "0"
$ws.Range("E114").Value = "2:06"
$ws.Range("F114").Value = "4:08"

$ws.Range("C115").Value = 1
$ws.Range("D115").NumberFormat = "@"
$ws.Range("D115").Value = "0"
$ws.Range("D115").NumberFormat = "0"
$ws.Range("E115").Value = "5:33"
$ws.Range("F115").Value = "6:52"

$ws.Range("C116").Value = 1
$ws.Range("D116").NumberFormat = "@"
$ws.Range("D116").Value = "0"
$ws.Range("D116").NumberFormat = "0"
$ws.Range("E116").Value = "7:41"
$ws.Range("F116").Value = "12:09"

$ws.Range("C117").Value = 2
$ws.Range("D117").NumberFormat = "@"
$ws.Range("D117").Value = "0"
$ws.Range("D117").NumberFormat = "0"
$ws.Range("E117").Value = "12:53"
$ws.Range("F117").Value = "16:03"

$ws.Range("C118").Value = 1
$ws.Range("D118").NumberFormat = "@"
$ws.Range("D118").Value = "0"
$ws.Range("D118").NumberFormat = "0"
$ws.Range("E118").Value = "19:58"
$ws.Range("F118").Value = "23:09"

$ws.Range("C119").Value = 1
$ws.Range("D119").NumberFormat = "@"
$ws.Range("D119").Value = "0"
$ws.Range("D119").NumberFormat = "0"
$ws.Range("E119").Value = "33:54"
$ws.Range("F119").Value = "35:32"

$ws.Range("C120").Value = 9
$ws.Range("D120").NumberFormat = "@"
$ws.Range("D120").Value = "0"
$ws.Range("D120").NumberFormat = "0"
$ws.Range("E120").Value = "36:01"
$ws.Range("F120").Value = "38:27"

$ws.Range("C121").Value = 3
$ws.Range("D121").NumberFormat = "@"
$ws.Range("D121").Value = "0"
$ws.Range("D121").NumberFormat = "0"
$ws.Range("E121").Value = "40:18"
$ws.Range("F121").Value = "42:28"

$ws.Range("C122").Value = 1
$ws.Range("D122").Value = 1
$ws.Range("E122").Value = "3:22"
$ws.Range("F122").Value = "9:23"

$ws.Range("C123").Value = 1
$ws.Range("D123").NumberFormat = "@"
$ws.Range("D123").Value = "0"
$ws.Range("D123").NumberFormat = "0"
$ws.Range("E123").Value = "11:01"
$ws.Range("F123").Value = "16:17"

$ws.Range("C124").Value = 1
$ws.Range("D124").Value = 2
$ws.Range("E124").Value = "17:00"
$ws.Range("F124").Value = "22:27"

$ws.Range("C125").Value = 4
$ws.Range("D125").NumberFormat = "@"
$ws.Range("D125").Value = "0"
$ws.Range("D125").NumberFormat = "0"
$ws.Range("E125").Value = "23:06"
$ws.Range("F125").Value = "27:04"

$ws.Range("C126").Value = 6
$ws.Range("D126").NumberFormat = "@"
$ws.Range("D126").Value = "0"
$ws.Range("D126").NumberFormat = "0"
$ws.Range("E126").Value = "32:06"
$ws.Range("F126").Value = "34:32"

$ws.Range("C127").Value = 1
$ws.Range("D127").NumberFormat = "@"
$ws.Range("D127").Value = "0"
$ws.Range("D127").NumberFormat = "0"
$ws.Range("E127").Value = "35:31"
$ws.Range("F127").Value = "37:48"

$ws.Range("C128").Value = 1
$ws.Range("D128").NumberFormat = "@"
$ws.Range("D128").Value = "0"
$ws.Range("D128").NumberFormat = "0"
$ws.Range("E128").Value = "38:26"
$ws.Range("F128").Value = "39:37"

$ws.Range("C129").Value = 2
$ws.Range("D129").NumberFormat = "@"
$ws.Range("D129").Value = "0"
$ws.Range("D129").NumberFormat = "0"
$ws.Range("E129").Value = "40:29"
$ws.Range("F129").Value = "42:50"
